# Parameters vs Performance of ML Model - add experiment run #19
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bring over the formatting (styles + row height) of the last data row (22)
# onto the new row (23), the same way a user would fill-down/copy the row.
$ws.Range("A22:E22").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's values
$ws.Range("A23").Value = 19
$ws.Range("B23").Value = 30000
$ws.Range("C23").Value = 0.0001
$ws.Range("D23").Value = 700
$ws.Range("E23").Value = $ws.Range("E22").Value()

# Make sure the row height matches row 22 (75pt)
$ws.Rows.Item(23).RowHeight = 75

# Update the selection to reflect where the user ended up after the edit
$ws.Range("D23").Select()
